$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-03-31 Monday"; new = "2025-04-01 Tuesday"},
    @{old = "673÷2=336, 1"; new = "531÷7=75, 6"},
    @{old = "862÷8=107, 6"; new = "772÷2=386, 0"},
    @{old = "205÷3=68, 1"; new = "719÷2=359, 1"},
    @{old = "138÷3=46, 0"; new = "443÷8=55, 3"},
    @{old = "735÷4=183, 3"; new = "847÷4=211, 3"},
    @{old = "269÷4=67, 1"; new = "105÷2=52, 1"},
    @{old = "786÷5=157, 1"; new = "847÷3=282, 1"},
    @{old = "188÷2=94, 0"; new = "555÷2=277, 1"},
    @{old = "120÷9=13, 3"; new = "614÷5=122, 4"},
    @{old = "199÷7=28, 3"; new = "855÷8=106, 7"},
    @{old = "201÷2=100, 1"; new = "834÷4=208, 2"},
    @{old = "356÷9=39, 5"; new = "200÷8=25, 0"},
    @{old = "912÷7=130, 2"; new = "959÷3=319, 2"},
    @{old = "856÷4=214, 0"; new = "361÷9=40, 1"},
    @{old = "827÷3=275, 2"; new = "520÷9=57, 7"},
    @{old = "259÷3=86, 1"; new = "714÷8=89, 2"},
    @{old = "188÷4=47, 0"; new = "695÷8=86, 7"},
    @{old = "140÷9=15, 5"; new = "144÷3=48, 0"},
    @{old = "467÷8=58, 3"; new = "144÷7=20, 4"},
    @{old = "845÷4=211, 1"; new = "469÷3=156, 1"},
    @{old = "430÷4=107, 2"; new = "273÷7=39, 0"},
    @{old = "134÷7=19, 1"; new = "573÷6=95, 3"},
    @{old = "197÷6=32, 5"; new = "166÷7=23, 5"},
    @{old = "171÷3=57, 0"; new = "351÷2=175, 1"},
    @{old = "327÷2=163, 1"; new = "886÷7=126, 4"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $r.new, 2)
}

$d.Save()
